$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Re-number the "smarthosting" bookmark (w:id 1 -> 0). Word assigns
#    bookmark ids sequentially starting at 0 when a bookmark is (re)added,
#    so deleting and re-adding it at the same range renumbers it.
# ---------------------------------------------------------------------
$bm = $d.Bookmarks("smarthosting")
$bmRange = $bm.Range
$bm.Delete() | Out-Null
$d.Bookmarks.Add("smarthosting", $bmRange) | Out-Null

# ---------------------------------------------------------------------
# 2. Remove the two HYPERLINK fields ("CryptoBridge" and "HitBTC") that
#    follow "...exchanges such as ". Field.Delete() removes the whole
#    field (begin/separate/end marks and displayed text) in one go.
# ---------------------------------------------------------------------
for ($i = $d.Fields.Count; $i -ge 1; $i--) {
    $f = $d.Fields.Item($i)
    if ($f.Code.Text -like "*crypto-bridge.org*") {
        $f.Delete() | Out-Null
    }
}
for ($i = $d.Fields.Count; $i -ge 1; $i--) {
    $f = $d.Fields.Item($i)
    if ($f.Code.Text -like "*hitbtc.com*") {
        $f.Delete() | Out-Null
    }
}

# ---------------------------------------------------------------------
# 3. Clean up the surrounding text so that
#      " Arguably the hardest part. Smart can be obtained from exchanges
#        such as <CryptoBridge>, <HitBTC>."
#    becomes
#      " Arguably the hardest part. Smart can be obtained from exchanges."
#    i.e. delete " such as" + the (now empty) ", " run that used to sit
#    between the two hyperlinks, leaving "exchanges" immediately followed
#    by the existing "." run.
# ---------------------------------------------------------------------
$rngFrom = $d.Content
$rngFrom.Find.Execute("exchanges") | Out-Null
$cutStart = $rngFrom.End

$rngTo = $d.Content
$rngTo.Start = $cutStart
$rngTo.End = $d.Content.End
$rngTo.Find.Execute(".") | Out-Null
$cutEnd = $rngTo.Start

$d.Range($cutStart, $cutEnd).Delete() | Out-Null
